$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 310, shifting existing rows 310:324 down to 311:325.
$ws.Rows("310:310").Insert()

# Populate the newly inserted row 310 with the new data record
# (same record shape/format as its neighbours, new date + volume).
$ws.Cells.Item(310, 1).Value = 9
$ws.Cells.Item(310, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(310, 3).Value = "Metropolitana"
$ws.Cells.Item(310, 4).Value = 44931
$ws.Cells.Item(310, 5).Value = 13
$ws.Cells.Item(310, 6).Value = 100112001
$ws.Cells.Item(310, 7).Value = "Berenjena"
$ws.Cells.Item(310, 8).Value = "Sin especificar"
$ws.Cells.Item(310, 9).Value = "Primera"
$ws.Cells.Item(310, 10).Value = 90
$ws.Cells.Item(310, 11).Value = 10000
$ws.Cells.Item(310, 12).Value = 12000
$ws.Cells.Item(310, 13).Value = 11000
$ws.Cells.Item(310, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(310, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(310, 16).Value = 220
$ws.Cells.Item(310, 17).Value = 50
$ws.Cells.Item(310, 18).Value = "Hortaliza"
